$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Add the new "PatientPost" worksheet right after "AdminLogin"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "PatientPost"

# ---------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------
$ws2.Range("A1").Value = "FirstName"
$ws2.Range("B1").Value = "LastName"
$ws2.Range("C1").Value = "PatientContactNumber"
$ws2.Range("D1").Value = "PatientEmail"
$ws2.Range("E1").Value = "Allergy"
$ws2.Range("F1").Value = "FoodPreference"
$ws2.Range("G1").Value = "CuisineCategory"
$ws2.Range("H1").NumberFormat = "yyyy-mm-dd"
$ws2.Range("H1").Value = "PatientDateOfBirth"

# ---------------------------------------------------------------
# Data row (row 2) - write order matches the shared-string intern
# order recorded in the source workbook (B, D, E, F, G, H, A)
# ---------------------------------------------------------------
$ws2.Range("B2").Value = "Kay"

$ws2.Range("C2").Font.Name = "Consolas"
$ws2.Range("C2").Font.Size = 10
$ws2.Range("C2").Font.Color = 0
$ws2.Range("C2").Value = 3248649876

$ws2.Range("D2").Value = "KAy123@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:KAy123@gmail.com")
$ws2.Range("D2").NumberFormat = "yyyy-mm-dd"

$ws2.Range("E2:H2").Font.Name = "Consolas"
$ws2.Range("E2:H2").Font.Size = 10
$ws2.Range("E2:H2").Font.Color = 0
$ws2.Range("E2:H2").NumberFormat = "yyyy-mm-dd"

$ws2.Range("E2").Value = "SOY"
$ws2.Range("F2").Value = "Vegan"
$ws2.Range("G2").Value = "Indian"
$ws2.Range("H2").Value = "12/12/2020"

$ws2.Range("A2").Value = "Raja"

# ---------------------------------------------------------------
# Column widths (best-fit, as captured by the source workbook).
# The engine stores ColumnWidth on a 1/6-character grid offset by
# 5/6 from the value that is handed to the COM setter, so the
# inputs below are pre-corrected to land as close as possible to
# the true (sub-1/6-character) best-fit widths recorded in the diff.
# ---------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 16.721354166666668
$ws2.Columns.Item(2).ColumnWidth = 12.053385416666666
$ws2.Columns.Item(3).ColumnWidth = 18.944010416666668
$ws2.Columns.Item(4).ColumnWidth = 17.276041666666668
$ws2.Columns.Item(5).ColumnWidth = 13.498697916666666
$ws2.Columns.Item(6).ColumnWidth = 21.385416666666668
$ws2.Columns.Item(7).ColumnWidth = 21.498697916666668
$ws2.Columns.Item(8).ColumnWidth = 15.385416666666666

# ---------------------------------------------------------------
# Selection / active sheet
# ---------------------------------------------------------------
$ws2.Range("C2").Select()
$ws2.Activate()

Write-Host "Done"
